$p = $ppt.ActivePresentation

# The deck currently applies the "Integral" theme (ppt/theme/theme1.xml) to the
# slide master, while the leftover default "Office Theme" only lives in
# ppt/theme/theme2.xml (used by the notes master). The edit swaps the two
# palettes back: the slide master goes back to the default Office palette.
#
# ColorScheme.Colors(1..12) maps 1:1 onto clrScheme's
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6,
# hlink, folHlink slots (RGB long values are R + G*256 + B*65536).

$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
